$d = $word.ActiveDocument

# 1) Merge the two runs that make up the "Comparative Ignorance..." heading
#    into a single run by re-finding the combined text and "replacing" it
#    with itself - Word's Find/Replace naturally collapses the matched
#    range into one run.
$d.Content.Find.Execute("Comparative Ignorance and the Ellsberg Paradox – Chow and Sarin", $true, $false, $false, $false, $false, $true, 1, $false, "Comparative Ignorance and the Ellsberg Paradox – Chow and Sarin", 2) | Out-Null

# 2) Replace the contents of the bookmarked (previously empty) paragraph
#    with the quoted excerpt, wrapped in curly quotes, while preserving
#    the _GoBack bookmark in place.
$bm = $d.Bookmarks("_GoBack")
$bmPara = $bm.Range.Paragraphs(1)
$bmParaRange = $d.Range($bmPara.Range.Start, $bmPara.Range.End)
$bmParaRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>‘</w:t></w:r><w:r><w:t>The key ﬁnding that emerges from our experiments is that the clear bet is priced higher than the vague bet under both comparative and non-comparative conditions. The comparison, however, enhances the difference in prices between clear and vague bets. In the absence of a direct comparison (non-comparative condition) this difference is smaller, but it does not disappear. This reduction in price differential between the clear and vague bets in the non-comparative condition is not evidence against ambiguity avoidance. Our results do not support the strong conclusion of Fox and Tversky (1995) that ambiguity aversion disappears in separate evaluations.</w:t></w:r><w:r><w:t>’</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

# 3) Merge the three runs that make up the "Is Luck on My Side..." heading
#    into a single run the same way as step 1.
$d.Content.Find.Execute("Is Luck on My Side? Optimism, Pessimism, and Ambiguity Aversion - Briony D. Pulford", $true, $false, $false, $false, $false, $true, 1, $false, "Is Luck on My Side? Optimism, Pessimism, and Ambiguity Aversion - Briony D. Pulford", 2) | Out-Null
